# "changement date gentleman nommay"
# The "Gentlemen de Nommay" race moved from "Dim 17 Octobre" (row 60) to
# "Dim 10 Octobre" (now before the Morschwiller cyclo-cross, so it becomes
# row 59). The Morschwiller race, previously row 59, shifts down to row 60
# and keeps its own date ("Sam 16 Octobre"), club and discipline, but its
# discipline cell now also reads "Route" in the edited workbook.
# Two unrelated small cleanups on the "Info" column (G) are also applied:
#  - G42 ("Magstatt le Bas") shortened from
#    "Annulé (nouvelles restrictions sanitaires depuis le 9 Août)" to "Annulé"
#  - G44 ("Prix de Fontaine") cleared (was "Anciennement Rougemont-le-Château")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59 becomes the Nommay (Gentlemen) entry, with its date moved up.
$ws.Cells.Item(59, 1).Value = "Dim 10 Octobre"
$ws.Cells.Item(59, 2).Value = "Gentlemen de Nommay (épreuve FFC ouverte aux FSGT - à confirmer)"
$ws.Cells.Item(59, 3).Value = "CCI Nommay"
$ws.Cells.Item(59, 4).Value = "Route"
$ws.Cells.Item(59, 5).Value = "nommay"

# Row 60 becomes the Morschwiller entry (was row 59), keeping its own date.
$ws.Cells.Item(60, 1).Value = "Sam 16 Octobre"
$ws.Cells.Item(60, 2).Value = "12<sup>e</sup> Cyclo-cross de l'Amitié et de la Solidarité à Morschwiller le bas  "
$ws.Cells.Item(60, 3).Value = "SOS Lutterbach"
$ws.Cells.Item(60, 4).Value = "Route"
$ws.Cells.Item(60, 5).Value = "morschwiller"

# Unrelated small text tweaks bundled in the same commit.
$ws.Cells.Item(42, 7).Value = "Annulé"
$ws.Cells.Item(44, 7).ClearContents()
